$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36 -- this shifts the existing rows 36..110 down
# to 37..111 (and their formatting/styles along with them), matching the
# diff's row-shift pattern.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Cells.Item(36, 1).Value = 6
$ws.Cells.Item(36, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(36, 3).Value = "Metropolitana"
$ws.Cells.Item(36, 4).Value = 44469
$ws.Cells.Item(36, 5).Value = 13
$ws.Cells.Item(36, 6).Value = 100112029
$ws.Cells.Item(36, 7).Value = "Orégano"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 32
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 9000
$ws.Cells.Item(36, 13).Value = 8469
$ws.Cells.Item(36, 14).Value = "$/docena de atados"
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 2823
$ws.Cells.Item(36, 17).Value = 3
$ws.Cells.Item(36, 18).Value = "Hortaliza"
